# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this game's individual play-by-play yardages to the
# season-long rush/pass logs (Offense Rush, Defense Rush, Offense Pass,
# Defense Pass).
# ---------------------------------------------------------------------------
$ydsSheet = $wb.Worksheets.Item("YDS")

$ydsSheet.Range("B2").Value = $ydsSheet.Range("B2").Value2 + " -3 5 3 5 2 4 1 3 16 1 3 16 13 -1 2 8 -1 5 8 8 0 2 3 3 1 8 -2 10 2 0 11"
$ydsSheet.Range("C2").Value = $ydsSheet.Range("C2").Value2 + " 6 1 0 3 1 -2 -1 -1 1 7 6 5 11 2 6 3 -3 15 3"
$ydsSheet.Range("B3").Value = $ydsSheet.Range("B3").Value2 + " 5 2 29 13 4 9 19 17 9 15 0 5 14 24 14 1 14 4 10 36 24 3 16"
$ydsSheet.Range("C3").Value = $ydsSheet.Range("C3").Value2 + " 7 8 14 6 -3 25 4 7 9 -5 20 28 24 8 13 37 30 7 5 5 21 5 26"

# ---------------------------------------------------------------------------
# OFF sheet: Home/Road offensive drive-chart totals after the Wild Card game.
# ---------------------------------------------------------------------------
$offSheet = $wb.Worksheets.Item("OFF")

$offSheet.Range("B2").Value = 13
$offSheet.Range("C2").Value = 520
$offSheet.Range("E2").Value = 21
$offSheet.Range("F2").Value = 168
$offSheet.Range("G2").Value = 121
$offSheet.Range("J2").Value = 70
$offSheet.Range("N2").Value = 38

$offSheet.Range("C3").Value = 268
$offSheet.Range("D3").Value = 13
$offSheet.Range("E3").Value = 66
$offSheet.Range("F3").Value = 199
$offSheet.Range("G3").Value = 69
$offSheet.Range("H3").Value = 54
$offSheet.Range("I3").Value = 114
$offSheet.Range("J3").Value = 123
$offSheet.Range("L3").Value = 497
$offSheet.Range("M3").Value = 326
$offSheet.Range("Q3").Value = 1099

# ---------------------------------------------------------------------------
# DEF sheet: Home/Road defensive drive-chart totals after the Wild Card game.
# ---------------------------------------------------------------------------
$defSheet = $wb.Worksheets.Item("DEF")

$defSheet.Range("C2").Value = 357
$defSheet.Range("D2").Value = 20
$defSheet.Range("E2").Value = 17
$defSheet.Range("F2").Value = 114
$defSheet.Range("G2").Value = 128
$defSheet.Range("H2").Value = 5
$defSheet.Range("I2").Value = 19
$defSheet.Range("J2").Value = 66
$defSheet.Range("N2").Value = 29
$defSheet.Range("O2").Value = 49
$defSheet.Range("P2").Value = 29

$defSheet.Range("B3").Value = 25
$defSheet.Range("C3").Value = 415
$defSheet.Range("E3").Value = 61
$defSheet.Range("F3").Value = 242
$defSheet.Range("G3").Value = 63
$defSheet.Range("H3").Value = 48
$defSheet.Range("I3").Value = 134
$defSheet.Range("J3").Value = 113
$defSheet.Range("L3").Value = 578
$defSheet.Range("M3").Value = 373
$defSheet.Range("Q3").Value = 995

# ---------------------------------------------------------------------------
# ST sheet: special-teams kickoff/punt counts & distance logs.
# ---------------------------------------------------------------------------
$stSheet = $wb.Worksheets.Item("ST")

$stSheet.Range("B2").Value = 186
$stSheet.Range("D2").Value = 110
$stSheet.Range("F2").Value = 270
$stSheet.Range("G2").Value = 257
$stSheet.Range("L2").Value = 83

$stSheet.Range("B3").Value = 108
$stSheet.Range("D3").Value = $stSheet.Range("D3").Value2 + " 53 46 49 49"

$stSheet.Range("B4").Value = $stSheet.Range("B4").Value2 + " 63 61 65"
$stSheet.Range("D4").Value = $stSheet.Range("D4").Value2 + " 13 3 8 0"

$stSheet.Range("B5").Value = $stSheet.Range("B5").Value2 + " 16 24 18"
$stSheet.Range("D5").Value = $stSheet.Range("D5").Value2 + " 1 0 24 8 10"

$stSheet.Range("B6").Value = $stSheet.Range("B6").Value2 + " 25 21 12"

# ---------------------------------------------------------------------------
# TURNS sheet: Road turnover counts.
# ---------------------------------------------------------------------------
$turnsSheet = $wb.Worksheets.Item("TURNS")

$turnsSheet.Range("D3").Value = 15
$turnsSheet.Range("E3").Value = 20

# ---------------------------------------------------------------------------
# PEN sheet: Holding penalty count.
# ---------------------------------------------------------------------------
$penSheet = $wb.Worksheets.Item("PEN")

$penSheet.Range("B3").Value = 31
